# edit.ps1
# Applies the 'Updated symbol list' commit changes to the cryptos worksheet.
# Column D values are numeric-looking but must remain TEXT (inline/shared string),
# matching the original workbook's inlineStr storage (e.g. '244.15', '0.03120').
# Columns B, C, E receive plain text updates (names, links, volume labels).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-looking values that must stay as TEXT (Price column D) ---
$priceUpdates = @{
    "D2" = '244.15'
    "D5" = '0.05610'
    "D6" = '3.365'
    "D7" = '6.382'
    "D8" = '0.8046'
    "D9" = '0.9854'
    "D10" = '0.01118'
    "D11" = '0.1419'
    "D12" = '0.07312'
    "D13" = '0.03120'
    "D14" = '0.03053'
    "D15" = '0.09288'
    "D16" = '3.573'
    "D17" = '0.001645'
    "D18" = '0.04706'
    "D19" = '0.006368'
    "D20" = '0.004993'
    "D21" = '0.001041'
    "D22" = '0.0001499'
    "D25" = '2.096'
    "D26" = '0.3260'
    "D40" = '0.03915'
    "D42" = '0.003397'
    "D43" = '0.1036'
    "D44" = '0.008493'
    "D45" = '0.00005941'
    "D47" = '0.0005493'
    "D48" = '0.6819'
    "D49" = '0.08686'
    "D50" = '0.00002098'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# --- Plain text values (Coin name, Link, Volume label columns) ---
$textUpdates = @{
    "B10" = 'One'
    "C10" = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    "E10" = '9OneONEBestin24h'
    "B11" = 'WazirX'
    "C11" = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    "E11" = '10WazirXWRX'
    "B12" = 'MandalaExchangeToken'
    "C12" = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    "E12" = '11MandalaExchangeTokenMDX'
    "B13" = 'LiechtensteinCryptoassetsExchange'
    "C13" = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    "E13" = '12LiechtensteinCryptoassetsExchangeLCX'
    "B14" = 'BitrueCoin'
    "C14" = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    "E14" = '13BitrueCoinBTR'
    "B15" = 'BitMartToken'
    "C15" = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    "E15" = '14BitMartTokenBMX'
    "B16" = 'MCDex'
    "C16" = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    "E16" = '15MCDexMCB'
    "B17" = 'BitForexToken'
    "C17" = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    "E17" = '16BitForexTokenBF'
    "B18" = 'CoinExToken'
    "C18" = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    "E18" = '17CoinExTokenCET'
    "B42" = 'CEJI'
    "C42" = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    "E42" = '41CEJICEJI'
    "B43" = 'BKEXToken'
    "C43" = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    "E43" = '42BKEXTokenBKK'
    "E47" = '46ACDXExchangeACXTWorstin24h'
    "E49" = '48BOLOBOLO'
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}
